$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 548, shifting existing rows 548:616 down to 549:617.
$ws.Rows.Item(548).Insert(-4121)

# Populate the constant columns (same for every data row in this block).
$ws.Cells.Item(548, 1).Value = 6
$ws.Cells.Item(548, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(548, 3).Value = "Metropolitana"
$ws.Cells.Item(548, 5).Value = 13
$ws.Cells.Item(548, 6).Value = 100112043
$ws.Cells.Item(548, 7).Value = "Pepino ensalada"
$ws.Cells.Item(548, 8).Value = "Sin especificar"
$ws.Cells.Item(548, 18).Value = "Hortaliza"

# Populate the new row's specific data values.
$ws.Cells.Item(548, 4).Value = 45077
$ws.Cells.Item(548, 9).Value = "Primera"
$ws.Cells.Item(548, 10).Value = 810
$ws.Cells.Item(548, 11).Value = 10000
$ws.Cells.Item(548, 12).Value = 11000
$ws.Cells.Item(548, 13).Value = 10568
$ws.Cells.Item(548, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(548, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(548, 16).Value = 176
$ws.Cells.Item(548, 17).Value = 60
